$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before row 1071, shifting the existing
# rows 1071:1127 down to 1074:1130.
$ws.Rows("1071:1073").Insert()

# Fill in the 3 new rows (1071:1073) with the new weekly entry
# (Fecha = 45041), mirroring the constant columns used throughout
# this block and the new M/N/O/P/S values from the diff.

# Row 1071 - Pintón
$ws.Range("A1071").Value = 8
$ws.Range("B1071").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1071").Value = "Coquimbo"
$ws.Range("D1071").Value = 45041
$ws.Range("E1071").Value = 4
$ws.Range("F1071").Value = "Fruta"
$ws.Range("G1071").Value = 100108
$ws.Range("H1071").Value = "Tropicales y subtropicales"
$ws.Range("I1071").Value = 100108006
$ws.Range("J1071").Value = "Plátano"
$ws.Range("K1071").Value = "Sin especificar"
$ws.Range("L1071").Value = "Pintón"
$ws.Range("M1071").Value = 120
$ws.Range("N1071").Value = 18000
$ws.Range("O1071").Value = 18000
$ws.Range("P1071").Value = 18000
$ws.Range("Q1071").Value = "$/caja 20 kilos"
$ws.Range("R1071").Value = "Ecuador"
$ws.Range("S1071").Value = 900
$ws.Range("T1071").Value = 20

# Row 1072 - Primera Maduro
$ws.Range("A1072").Value = 8
$ws.Range("B1072").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1072").Value = "Coquimbo"
$ws.Range("D1072").Value = 45041
$ws.Range("E1072").Value = 4
$ws.Range("F1072").Value = "Fruta"
$ws.Range("G1072").Value = 100108
$ws.Range("H1072").Value = "Tropicales y subtropicales"
$ws.Range("I1072").Value = 100108006
$ws.Range("J1072").Value = "Plátano"
$ws.Range("K1072").Value = "Sin especificar"
$ws.Range("L1072").Value = "Primera Maduro"
$ws.Range("M1072").Value = 120
$ws.Range("N1072").Value = 20000
$ws.Range("O1072").Value = 20000
$ws.Range("P1072").Value = 20000
$ws.Range("Q1072").Value = "$/caja 20 kilos"
$ws.Range("R1072").Value = "Ecuador"
$ws.Range("S1072").Value = 1000
$ws.Range("T1072").Value = 20

# Row 1073 - Primera Pintón
$ws.Range("A1073").Value = 8
$ws.Range("B1073").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1073").Value = "Coquimbo"
$ws.Range("D1073").Value = 45041
$ws.Range("E1073").Value = 4
$ws.Range("F1073").Value = "Fruta"
$ws.Range("G1073").Value = 100108
$ws.Range("H1073").Value = "Tropicales y subtropicales"
$ws.Range("I1073").Value = 100108006
$ws.Range("J1073").Value = "Plátano"
$ws.Range("K1073").Value = "Sin especificar"
$ws.Range("L1073").Value = "Primera Pintón"
$ws.Range("M1073").Value = 120
$ws.Range("N1073").Value = 21000
$ws.Range("O1073").Value = 21000
$ws.Range("P1073").Value = 21000
$ws.Range("Q1073").Value = "$/caja 20 kilos"
$ws.Range("R1073").Value = "Ecuador"
$ws.Range("S1073").Value = 1050
$ws.Range("T1073").Value = 20
